$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "Lime" translation row for D column, mirroring existing Lily/Shina entries in C,
# plus duplicate the MP_SET_LOOP command into D8 as in the source diff.
$ws.Range("D2").Value = "Lily"
$ws.Range("D3").Value = "Shina"
$ws.Range("D4").Value = "Lime"
$ws.Range("D8").Value = "MP_SET_LOOP 20 on"
